$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Before state (rows 4-8): AUGMENTIN, BRUFEN, VOLTAREN, WATER FOR INJECTION,
# سرنجات 3 سم, followed by a totals row and a footer row.
#
# After state: 4 new medicine rows (FLACORT, GLUCOPHAGE XR, IVYROSPAN,
# MAXILASE) are inserted between BRUFEN (row 5) and VOLTAREN (old row 6),
# pushing VOLTAREN/WATER/سرنجات and the totals/footer rows down by 4.
# ---------------------------------------------------------------------------

# Insert 4 blank rows at row 6 (pushes VOLTAREN.. down to rows 10-12,
# the totals row down to 13 and the footer row down to 14).
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(6).Insert()
}

# Give each newly inserted row the same formatting (styles + borders) as
# an existing data row. After the 4 inserts above, the original row 8
# ("سرنجات 3 سم") has been pushed down to row 12, so copy the format
# from there (copying from a still-blank row would just propagate the
# border-less blank-row style).
for ($r = 6; $r -le 9; $r++) {
    $ws.Range("A12:N12").Copy()
    $ws.Range("A" + $r + ":N" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Row heights match the auto-fit heights Excel computed for this content
# (see target sheet XML).
$ws.Rows.Item(6).RowHeight = 24.75
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 24.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 24.75
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 17.25

# Re-create the merged cells (B:G / H:K / L:M) on each new row.
foreach ($r in 6,7,8,9) {
    $ws.Range("B" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
}

# Fill in the 4 new medicine rows.
$ws.Cells.Item(6, 1).Value = 3
$ws.Cells.Item(6, 2).Value = "FLACORT 30MG 20 TAB"
$ws.Cells.Item(6, 8).Value = "1:0"
$ws.Cells.Item(6, 12).Value = 82
$ws.Cells.Item(6, 14).Value = "0:2"

$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "GLUCOPHAGE XR 1000 MG 30 TABS."
$ws.Cells.Item(7, 8).Value = "5:2"
$ws.Cells.Item(7, 12).Value = 42
$ws.Cells.Item(7, 14).Value = "0:0"

$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "IVYROSPAN SYRUP 100 ML"
$ws.Cells.Item(8, 8).Value = "0:0"
$ws.Cells.Item(8, 12).Value = 55
$ws.Cells.Item(8, 14).Value = "1:0"

$ws.Cells.Item(9, 1).Value = 6
$ws.Cells.Item(9, 2).Value = "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML"
$ws.Cells.Item(9, 8).Value = "4:0"
$ws.Cells.Item(9, 12).Value = 57
$ws.Cells.Item(9, 14).Value = "1:0"

# Renumber the "م" (sequence) column for the rows that got pushed down.
$ws.Cells.Item(10, 1).Value = 7   # VOLTAREN 75MG/3ML 3 AMP.
$ws.Cells.Item(11, 1).Value = 8   # WATER FOR INJECTION AMP. 5 ML
$ws.Cells.Item(12, 1).Value = 9   # سرنجات 3 سم

# Update the grand total (sum of column L) on the totals row, now at row 13.
$ws.Cells.Item(13, 11).Value = 403

Write-Host "done"
